# NegronSierra_ProjectRubric.xlsx edit
# "Moved the parsing and Game Level data over to the renderer for efficiency.
#  Partial implementation of level swap."
#
# Net effect on Sheet1:
#  - The comment that used to sit in E9 ("...some submeshes aren't drawing, it
#    seems") is replaced by a new, fuller comment in D9.
#  - Short "d" markers are added in E14, E15 and E21.
#  - The selected cell moves from C23 to F13.
#
# NOTE on ordering: the shared-strings table is rebuilt from scratch on save,
# assigning ids as new unique strings are first encountered/created. To land
# on the same ids as the target workbook (short "d" marker before the long
# "Partially complete..." comment), E9 must be cleared (freeing its old
# string) and the "d" cells must be written *before* the new D9 text is
# written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old comment from E9 first.
$ws.Range("E9").ClearContents()

# Add the short "d" markers.
$ws.Range("E14").Value = "d"
$ws.Range("E15").Value = "d"
$ws.Range("E21").Value = "d"

# Write the new, expanded comment into D9.
$ws.Range("D9").Value = "Partially complete, some submeshes aren't drawing, it seems"

# Update the active selection to match the saved workbook state.
$ws.Range("F13").Select()
